$wb = $excel.ActiveWorkbook
$wsSites = $wb.Worksheets.Item("Sites")
$wsIndividuals = $wb.Worksheets.Item("Individuals")

# Insert a new "Comments" column before the "Count Code" column (G) on the
# Sites sheet, shifting Count Code/Count right by one.
$wsSites.Columns("G").Insert()
$wsSites.Columns("G").ColumnWidth = $wsSites.Columns("F").ColumnWidth
$wsSites.Range("G3").Value = "Comments"

# Sites becomes the active sheet/tab, selection resting on the new header cell.
$wsIndividuals.Range("F7").Select()
$wsSites.Activate()
$wsSites.Range("G3").Select()
